$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test2")

$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 26
$ws.Range("D4").Value = 50
